$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'327.43"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'-1.05%"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'43.77"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'5.36%"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = "'5.467"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'0.08079"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'-3.75%"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'8.652"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'-1.93%"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'4.293"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'-4.22%"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'1.882"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'-6.07%"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'-8.21%"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'0.9383"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'1.47%"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.1187"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'-7.24%"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = "'0.1896"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'-3.21%"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'0.09564"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'2.54%"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'0.04105"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'3.36%"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'0.60%"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'0.001281"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'-1.97%"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'0.005932"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'-2.96%"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'3.577"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'4.47%"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'-0.05%"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = "'8.643"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'-3.35%"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'-1.20%"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'0.2496"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'-0.64%"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'0.04366"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'-1.00%"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'0.001236"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'-0.78%"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'0.004309"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'-1.35%"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = "'0.0001235"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'3.65%"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.Value = "'0.0004005"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'0.25%"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'0.02656"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'-6.43%"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.05417"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'-1.82%"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.007614"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'-4.38%"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'0.01002"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'11.74%"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.1390"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'-3.18%"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.002094"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'0.54%"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'0.009907"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'-10.99%"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'0.00006881"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'-1.10%"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.00000000753"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'0.26%"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'0.003563"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'8.00%"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = "'0.002279"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'-0.08%"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'0.00002108"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'0.26%"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'0.0002008"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'0.26%"
$cell.Style = "Normal"
